$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 314 values (D314, F314)
$ws.Range("D314").Value = 8.0829
$ws.Range("F314").Value = 8.0769

# Copy the date-column formatting (style) from row 314 down to the new rows 315:317
$ws.Range("A314").Copy()
$ws.Range("A315:A317").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add new row 315 (2023-09-01)
$ws.Range("A315").Value = 45170.33333333334
$ws.Range("B315").Value = "FX_IDC:USDMOP"
$ws.Range("C315").Value = 8.0771
$ws.Range("D315").Value = 8.0823
$ws.Range("E315").Value = 8.048400000000001
$ws.Range("F315").Value = 8.0654
$ws.Range("G315").Value = 0

# Add new row 316 (2023-10-02)
$ws.Range("A316").Value = 45201.375
$ws.Range("B316").Value = "FX_IDC:USDMOP"
$ws.Range("C316").Value = 8.0654
$ws.Range("D316").Value = 8.0684
$ws.Range("E316").Value = 8.0509
$ws.Range("F316").Value = 8.058400000000001
$ws.Range("G316").Value = 0

# Add new row 317 (2023-11-01)
$ws.Range("A317").Value = 45231.375
$ws.Range("B317").Value = "FX_IDC:USDMOP"
$ws.Range("C317").Value = 8.0585
$ws.Range("D317").Value = 8.0616
$ws.Range("E317").Value = 8.0381
$ws.Range("F317").Value = 8.043100000000001
$ws.Range("G317").Value = 0
